$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Planning": mark Bloc 5 column with an X for "descente de code" (row
# 36) and "Tests de l'application" (row 38).
# ---------------------------------------------------------------------------
$planning = $wb.Worksheets.Item("Planning")
$planning.Range("O36").Value2 = "X"
$planning.Range("O38").Value2 = "X"

# ---------------------------------------------------------------------------
# Sheet "Journal de travail": add a new day entry (03.03.2024) with four
# tasks, then shift the existing "Total" / legend block down to make room.
# ---------------------------------------------------------------------------
$journal = $wb.Worksheets.Item("Journal de travail")

# Insert six new rows right after the existing blank row 43 (i.e. before the
# old row 44), matching the formatting of row 43 so the new block keeps the
# same borders/number formats as the rest of the table.
$journal.Range("A44:A49").EntireRow.Insert()
$journal.Range("A43:C43").Copy()
$journal.Range("A44:C49").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 44's "hours" cell keeps the divider row's original number format
# (it was pushed down to row 50 by the insert above).
$journal.Range("C50").Copy()
$journal.Range("C44").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# The divider row (old row 44, now row 50) only keeps its first cell.
$journal.Range("B50:C50").Clear()

# Fill in the new day block.
$journal.Range("A43").Value2 = 45354
$journal.Range("B43").Value2 = "correction de bugs et nettoyage du code"
$journal.Range("C43").Value2 = 1

$journal.Range("B44").Value2 = "descente de code"
$journal.Range("C44").Value2 = 2

$journal.Range("B45").Value2 = "Tests fonctionnels"
$journal.Range("C45").Value2 = 2

$journal.Range("B46").Value2 = "JS et PHP doc"
$journal.Range("C46").Value2 = 1

# Rows 47-49 stay empty (already copied with the right formatting).

# Fix up the Total formula, which the row insertion widened to C8:C50 (it
# should stop at the new last data row, C49).
$journal.Range("C51").Formula = "=SUM(C8:C49)"

# Restore the selections recorded in the target workbook.
$planning.Select()
$planning.Range("O33").Select()
$journal.Select()
$journal.Range("B47").Select()
